# Generate Report for Handoff
#
# "b.md" has now been handed off for localization (zh-cn and de-de), so the
# status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" everywhere it is tracked, new handoff-xlf filenames/timestamps
# are recorded, and the per-language tables gain an error detail message
# about the handback file being stale.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is "b.md" ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 00:45:18"

# ---- zh-cn sheet: row 3 is "b.md" ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 00:45:13"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98ef8b00ee7e494e5cc6a8e402fbf3d9fefca25d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a44efa4f267435f83d08552e2948e166a7d70797/e2e/b.md."
$zhcn.Range("P1").ColumnWidth = 39.14

# ---- de-de sheet: row 3 is "b.md" ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 00:45:18"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98ef8b00ee7e494e5cc6a8e402fbf3d9fefca25d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a44efa4f267435f83d08552e2948e166a7d70797/e2e/b.md."
$dede.Range("P1").ColumnWidth = 39.14
